$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Switzerland" row for year 2023 (pushes 2025..2050 + the whole
#     "Vaud" formula block below it down by one row) ---
$ws.Rows.Item(35).Insert()
$ws.Range("A35").Value = 2023
$ws.Range("B35").Value = $ws.Range("B34").Value2
$ws.Range("C35").Value = $ws.Range("C34").Value2
$ws.Range("D35").Value = $ws.Range("D34").Value2
$ws.Range("E35").Value = $ws.Range("E34").Value2

# --- Insert a new "Vaud" row for year 2023 right after the (now shifted) 2022
#     row, mirroring the same 0.09*<Switzerland row> relationship as its peers ---
$ws.Rows.Item(75).Insert()
$ws.Range("A75").Value = 2023
$ws.Range("B75").Value = $ws.Range("B74").Value2
$ws.Range("C75").Formula = "=0.09*C35"
$ws.Range("D75").Formula = "=0.09*D35"
$ws.Range("E75").Formula = "=0.09*E35"

# --- Widen the data columns that now hold the (wider) inserted content ---
$ws.Columns.Item(3).ColumnWidth = 27.583333333333332
$ws.Columns.Item(4).ColumnWidth = 21.583333333333332
$ws.Columns.Item(5).ColumnWidth = 24.166666666666668

# --- Leave the selection on the newly-inserted "Vaud" 2023 row, matching the
#     saved view state ---
$ws.Range("A74:E75").Select()
